$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1818817.6
$ws.Range("I17").Value = 663
$ws.Range("J17").Value = 1852487.1
$ws.Range("K17").Value = 1989
$ws.Range("L17").Value = 5557461.300000001
$ws.Range("M17").Value = -1821
$ws.Range("N17").Value = -5557797.300000001
# Row 98
$ws.Range("H98").Value = 354.83334
$ws.Range("I98").Value = 369.82352
$ws.Range("J98").Value = 100
$ws.Range("K98").Value = 369.82352
$ws.Range("L98").Value = 100
$ws.Range("M98").Value = 1128.17648
# Row 116
$ws.Range("H116").Value = 3316
$ws.Range("I116").Value = 1461
$ws.Range("J116").Value = 3934.3333
$ws.Range("K116").Value = 1461
$ws.Range("L116").Value = 3934.3333
$ws.Range("M116").Value = 1981
$ws.Range("N116").Value = -10818.3333
# Row 122
$ws.Range("H122").Value = 354.83334
$ws.Range("I122").Value = 369.82352
$ws.Range("J122").Value = 100
$ws.Range("K122").Value = 1109.47056
$ws.Range("L122").Value = 300
$ws.Range("M122").Value = 1340.52944
# Row 129
$ws.Range("H129").Value = 185986.83
$ws.Range("I129").Value = 339.8
$ws.Range("J129").Value = 204930.4
$ws.Range("K129").Value = 1019.4
$ws.Range("L129").Value = 614791.2
$ws.Range("M129").Value = 3980.6
$ws.Range("N129").Value = -624791.2
# Row 137
$ws.Range("H137").Value = 171850.17
$ws.Range("I137").Value = 11500
$ws.Range("J137").Value = 252025.25
$ws.Range("K137").Value = 34500
$ws.Range("L137").Value = 756075.75
$ws.Range("M137").Value = -31950
$ws.Range("N137").Value = -761175.75

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2656.3333
$ws.Range("I45").Value = 3742.8572
$ws.Range("J45").Value = 1964.909
$ws.Range("K45").Value = 3742.8572
$ws.Range("L45").Value = 1964.909
$ws.Range("M45").Value = -3365.8572
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 2442.6
$ws.Range("I107").Value = 900
$ws.Range("J107").Value = 4756.5
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 4756.5
$ws.Range("M107").Value = 1020

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 13099.846
$ws.Range("I31").Value = 15418.857
$ws.Range("J31").Value = 3360
$ws.Range("K31").Value = 15418.857
$ws.Range("L31").Value = 3360
$ws.Range("M31").Value = -15123.857
$ws.Range("N31").Value = -3950
# Row 34
$ws.Range("H34").Value = 13099.846
$ws.Range("I34").Value = 15418.857
$ws.Range("J34").Value = 3360
$ws.Range("K34").Value = 15418.857
$ws.Range("L34").Value = 3360
$ws.Range("M34").Value = -15216.857
$ws.Range("N34").Value = -3764

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 4766.0625
$ws.Range("I3").Value = 1825.5555
$ws.Range("J3").Value = 8546.714
$ws.Range("K3").Value = 5476.666499999999
$ws.Range("L3").Value = 25640.142
$ws.Range("M3").Value = -5364.666499999999
$ws.Range("N3").Value = -25864.142
# Row 86
$ws.Range("H86").Value = 500000000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 500000000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 1500000000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -1500002372
# Row 89
$ws.Range("H89").Value = 500000000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 500000000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 4500000000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -4500011856
# Row 92
$ws.Range("H92").Value = 599.2222
$ws.Range("I92").Value = 628.5714
$ws.Range("J92").Value = 496.5
$ws.Range("K92").Value = 1885.7142
$ws.Range("L92").Value = 1489.5
$ws.Range("M92").Value = -637.7142000000001
# Row 93
$ws.Range("H93").Value = 7780.7144
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 7780.7144
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 23342.1432
$ws.Range("N93").Value = -27086.1432
# Row 98
$ws.Range("H98").Value = 358.7857
$ws.Range("I98").Value = 256.75
$ws.Range("J98").Value = 399.6
$ws.Range("K98").Value = 770.25
$ws.Range("L98").Value = 1198.8
$ws.Range("M98").Value = 727.75
$ws.Range("N98").Value = -4194.8
# Row 113
$ws.Range("H113").Value = 9075.75
$ws.Range("I113").Value = 33767
$ws.Range("J113").Value = 845.3333
$ws.Range("K113").Value = 101301
$ws.Range("L113").Value = 2535.9999
$ws.Range("M113").Value = -99131
$ws.Range("N113").Value = -6875.9999
# Row 131
$ws.Range("H131").Value = 746.1900000000001
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 746.1900000000001
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2238.57
$ws.Range("N131").Value = -12318.57

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
# Row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 57
$ws.Range("H57").Value = 28695
$ws.Range("I57").Value = 27400
$ws.Range("J57").Value = 29990
$ws.Range("K57").Value = 27400
$ws.Range("L57").Value = 29990
$ws.Range("M57").Value = -26580
$ws.Range("N57").Value = -31630
# Row 80
$ws.Range("H80").Value = 3117.2334
$ws.Range("I80").Value = 2739.9
$ws.Range("J80").Value = 3305.9
$ws.Range("K80").Value = 2739.9
$ws.Range("L80").Value = 3305.9
$ws.Range("M80").Value = -1741.9
$ws.Range("N80").Value = -5301.9
# Row 83
$ws.Range("H83").Value = 3117.2334
$ws.Range("I83").Value = 2739.9
$ws.Range("J83").Value = 3305.9
$ws.Range("K83").Value = 13699.5
$ws.Range("L83").Value = 16529.5
$ws.Range("M83").Value = -8707.5
$ws.Range("N83").Value = -26513.5
# Row 102
$ws.Range("H102").Value = 26318212
$ws.Range("I102").Value = 35717224
$ws.Range("J102").Value = 976.8
$ws.Range("K102").Value = 35717224
$ws.Range("L102").Value = 976.8
$ws.Range("M102").Value = -35715602
# Row 126
$ws.Range("H126").Value = 3710.8462
$ws.Range("I126").Value = 2843.1924
$ws.Range("J126").Value = 5446.154
$ws.Range("K126").Value = 8529.5772
$ws.Range("L126").Value = 16338.462
$ws.Range("M126").Value = -6059.5772
$ws.Range("N126").Value = -21278.462

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2933.4443
$ws.Range("I22").Value = 1985.8572
$ws.Range("J22").Value = 6250
$ws.Range("K22").Value = 1985.8572
$ws.Range("L22").Value = 6250
$ws.Range("M22").Value = -1690.8572
$ws.Range("N22").Value = -6840
# Row 27
$ws.Range("H27").Value = 2933.4443
$ws.Range("I27").Value = 1985.8572
$ws.Range("J27").Value = 6250
$ws.Range("K27").Value = 1985.8572
$ws.Range("L27").Value = 6250
$ws.Range("M27").Value = -1878.8572
$ws.Range("N27").Value = -6464

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 17242774
$ws.Range("I136").Value = 25001222
$ws.Range("J136").Value = 1776.2778
$ws.Range("K136").Value = 75003666
$ws.Range("L136").Value = 5328.8334
$ws.Range("M136").Value = -75001116
$ws.Range("N136").Value = -10428.8334
